$d = $word.ActiveDocument

# --- 1. Remove the _GoBack bookmark from its current location (end of doc) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Insert a period right after "...age gap (e.g. 100)" ---
# Remember the boundary between the "Maximum age-gap..." run and the
# " If you do not want..." run so we can restore it as a separate run
# after the engine's adjacent-same-format-run merge pass.
$rBoundary = $d.Content
$rBoundary.Find.Execute("Maximum age-gap for maternal half-siblings.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundaryB = $rBoundary.End

$rTarget = $d.Content
$rTarget.Find.Execute("use an arbitrarily large age gap (e.g. 100)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos = $rTarget.End

$ip = $d.Range($insertPos, $insertPos)
$ip.InsertAfter(".")

# --- 3. Undo the accidental merges the insertion caused ---
# a) restore " If you do not want ... (e.g. 100)" as its own run
$origRun = $d.Range($boundaryB, $insertPos)
$origRun.Bold = $true
$origRun.Bold = $false

# b) make sure the new "." stays its own run (not merged into the run
#    that follows it)
$dotRange = $d.Range($insertPos, $insertPos + 1)
$dotRange.Bold = $true
$dotRange.Bold = $false

# --- 4. Re-add the _GoBack bookmark right after the new period ---
$bmRange = $d.Range($insertPos + 1, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
